# Weekly fruit/vegetable price update: shift existing Brócoli price rows down by one
# (154..190 <- 153..189) and populate row 153 with a new weekly observation.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 153
$ws.Cells.Item(153, 4).Value = 44511
$ws.Cells.Item(153, 9).Value = "Primera"
$ws.Cells.Item(153, 10).Value = 300
$ws.Cells.Item(153, 11).Value = 650
$ws.Cells.Item(153, 12).Value = 700
$ws.Cells.Item(153, 13).Value = 675
$ws.Cells.Item(153, 15).Value = "Región del Maule"
$ws.Cells.Item(153, 16).Value = 675

# Row 154
$ws.Cells.Item(154, 4).Value = 44306
$ws.Cells.Item(154, 9).Value = "Primera"
$ws.Cells.Item(154, 10).Value = 300
$ws.Cells.Item(154, 11).Value = 700
$ws.Cells.Item(154, 12).Value = 750
$ws.Cells.Item(154, 13).Value = 725
$ws.Cells.Item(154, 15).Value = "Región del Maule"
$ws.Cells.Item(154, 16).Value = 725

# Row 155
$ws.Cells.Item(155, 4).Value = 44469
$ws.Cells.Item(155, 9).Value = "Primera"
$ws.Cells.Item(155, 10).Value = 300
$ws.Cells.Item(155, 11).Value = 700
$ws.Cells.Item(155, 12).Value = 750
$ws.Cells.Item(155, 13).Value = 725
$ws.Cells.Item(155, 15).Value = "Región del Maule"
$ws.Cells.Item(155, 16).Value = 725

# Row 156
$ws.Cells.Item(156, 4).Value = 44407
$ws.Cells.Item(156, 9).Value = "Primera"
$ws.Cells.Item(156, 10).Value = 160
$ws.Cells.Item(156, 11).Value = 550
$ws.Cells.Item(156, 12).Value = 600
$ws.Cells.Item(156, 13).Value = 575
$ws.Cells.Item(156, 15).Value = "Provincia de Diguillín"
$ws.Cells.Item(156, 16).Value = 575

# Row 157
$ws.Cells.Item(157, 4).Value = 44407
$ws.Cells.Item(157, 9).Value = "Segunda"
$ws.Cells.Item(157, 10).Value = 160
$ws.Cells.Item(157, 11).Value = 400
$ws.Cells.Item(157, 12).Value = 450
$ws.Cells.Item(157, 13).Value = 425
$ws.Cells.Item(157, 15).Value = "Provincia de Diguillín"
$ws.Cells.Item(157, 16).Value = 425

# Row 158
$ws.Cells.Item(158, 4).Value = 44186
$ws.Cells.Item(158, 9).Value = "Primera"
$ws.Cells.Item(158, 10).Value = 240
$ws.Cells.Item(158, 11).Value = 700
$ws.Cells.Item(158, 12).Value = 750
$ws.Cells.Item(158, 13).Value = 725
$ws.Cells.Item(158, 15).Value = "Provincia de Diguillín"
$ws.Cells.Item(158, 16).Value = 725

# Row 159
$ws.Cells.Item(159, 4).Value = 44504
$ws.Cells.Item(159, 9).Value = "Primera"
$ws.Cells.Item(159, 10).Value = 300
$ws.Cells.Item(159, 11).Value = 650
$ws.Cells.Item(159, 12).Value = 700
$ws.Cells.Item(159, 13).Value = 675
$ws.Cells.Item(159, 15).Value = "Región del Maule"
$ws.Cells.Item(159, 16).Value = 675

# Row 160
$ws.Cells.Item(160, 4).Value = 44384
$ws.Cells.Item(160, 9).Value = "Primera"
$ws.Cells.Item(160, 10).Value = 160
$ws.Cells.Item(160, 11).Value = 600
$ws.Cells.Item(160, 12).Value = 650
$ws.Cells.Item(160, 13).Value = 625
$ws.Cells.Item(160, 15).Value = "Región del Maule"
$ws.Cells.Item(160, 16).Value = 625

# Row 161
$ws.Cells.Item(161, 4).Value = 44246
$ws.Cells.Item(161, 9).Value = "Primera"
$ws.Cells.Item(161, 10).Value = 300
$ws.Cells.Item(161, 11).Value = 700
$ws.Cells.Item(161, 12).Value = 750
$ws.Cells.Item(161, 13).Value = 725
$ws.Cells.Item(161, 15).Value = "Región del Maule"
$ws.Cells.Item(161, 16).Value = 725

# Row 162
$ws.Cells.Item(162, 4).Value = 44505
$ws.Cells.Item(162, 9).Value = "Primera"
$ws.Cells.Item(162, 10).Value = 400
$ws.Cells.Item(162, 11).Value = 650
$ws.Cells.Item(162, 12).Value = 700
$ws.Cells.Item(162, 13).Value = 675
$ws.Cells.Item(162, 15).Value = "Región del Maule"
$ws.Cells.Item(162, 16).Value = 675

# Row 163
$ws.Cells.Item(163, 4).Value = 44487
$ws.Cells.Item(163, 9).Value = "Primera"
$ws.Cells.Item(163, 10).Value = 100
$ws.Cells.Item(163, 11).Value = 750
$ws.Cells.Item(163, 12).Value = 800
$ws.Cells.Item(163, 13).Value = 775
$ws.Cells.Item(163, 15).Value = "Región del Maule"
$ws.Cells.Item(163, 16).Value = 775

# Row 164
$ws.Cells.Item(164, 4).Value = 44425
$ws.Cells.Item(164, 9).Value = "Primera"
$ws.Cells.Item(164, 10).Value = 300
$ws.Cells.Item(164, 11).Value = 700
$ws.Cells.Item(164, 12).Value = 750
$ws.Cells.Item(164, 13).Value = 725
$ws.Cells.Item(164, 15).Value = "Región del Maule"
$ws.Cells.Item(164, 16).Value = 725

# Row 165
$ws.Cells.Item(165, 4).Value = 44343
$ws.Cells.Item(165, 9).Value = "Primera"
$ws.Cells.Item(165, 10).Value = 120
$ws.Cells.Item(165, 11).Value = 700
$ws.Cells.Item(165, 12).Value = 750
$ws.Cells.Item(165, 13).Value = 725
$ws.Cells.Item(165, 15).Value = "Región del Maule"
$ws.Cells.Item(165, 16).Value = 725

# Row 166
$ws.Cells.Item(166, 4).Value = 44370
$ws.Cells.Item(166, 9).Value = "Primera"
$ws.Cells.Item(166, 10).Value = 300
$ws.Cells.Item(166, 11).Value = 600
$ws.Cells.Item(166, 12).Value = 650
$ws.Cells.Item(166, 13).Value = 625
$ws.Cells.Item(166, 15).Value = "Provincia de Diguillín"
$ws.Cells.Item(166, 16).Value = 625

# Row 167
$ws.Cells.Item(167, 4).Value = 44370
$ws.Cells.Item(167, 9).Value = "Primera"
$ws.Cells.Item(167, 10).Value = 160
$ws.Cells.Item(167, 11).Value = 500
$ws.Cells.Item(167, 12).Value = 550
$ws.Cells.Item(167, 13).Value = 525
$ws.Cells.Item(167, 15).Value = "Región del Maule"
$ws.Cells.Item(167, 16).Value = 525

# Row 168
$ws.Cells.Item(168, 4).Value = 44449
$ws.Cells.Item(168, 9).Value = "Primera"
$ws.Cells.Item(168, 10).Value = 160
$ws.Cells.Item(168, 11).Value = 700
$ws.Cells.Item(168, 12).Value = 750
$ws.Cells.Item(168, 13).Value = 725
$ws.Cells.Item(168, 15).Value = "Región del Maule"
$ws.Cells.Item(168, 16).Value = 725

# Row 169
$ws.Cells.Item(169, 4).Value = 44168
$ws.Cells.Item(169, 9).Value = "Primera"
$ws.Cells.Item(169, 10).Value = 120
$ws.Cells.Item(169, 11).Value = 700
$ws.Cells.Item(169, 12).Value = 750
$ws.Cells.Item(169, 13).Value = 725
$ws.Cells.Item(169, 15).Value = "Región del Maule"
$ws.Cells.Item(169, 16).Value = 725

# Row 170
$ws.Cells.Item(170, 4).Value = 44175
$ws.Cells.Item(170, 9).Value = "Primera"
$ws.Cells.Item(170, 10).Value = 300
$ws.Cells.Item(170, 11).Value = 700
$ws.Cells.Item(170, 12).Value = 750
$ws.Cells.Item(170, 13).Value = 725
$ws.Cells.Item(170, 15).Value = "Región del Maule"
$ws.Cells.Item(170, 16).Value = 725

# Row 171
$ws.Cells.Item(171, 4).Value = 44392
$ws.Cells.Item(171, 9).Value = "Primera"
$ws.Cells.Item(171, 10).Value = 300
$ws.Cells.Item(171, 11).Value = 700
$ws.Cells.Item(171, 12).Value = 750
$ws.Cells.Item(171, 13).Value = 725
$ws.Cells.Item(171, 15).Value = "Provincia de Diguillín"
$ws.Cells.Item(171, 16).Value = 725

# Row 172
$ws.Cells.Item(172, 4).Value = 44286
$ws.Cells.Item(172, 9).Value = "Primera"
$ws.Cells.Item(172, 10).Value = 300
$ws.Cells.Item(172, 11).Value = 700
$ws.Cells.Item(172, 12).Value = 750
$ws.Cells.Item(172, 13).Value = 725
$ws.Cells.Item(172, 15).Value = "Provincia de Diguillín"
$ws.Cells.Item(172, 16).Value = 725

# Row 173
$ws.Cells.Item(173, 4).Value = 44473
$ws.Cells.Item(173, 9).Value = "Primera"
$ws.Cells.Item(173, 10).Value = 300
$ws.Cells.Item(173, 11).Value = 600
$ws.Cells.Item(173, 12).Value = 650
$ws.Cells.Item(173, 13).Value = 625
$ws.Cells.Item(173, 15).Value = "Provincia de Diguillín"
$ws.Cells.Item(173, 16).Value = 625

# Row 174
$ws.Cells.Item(174, 4).Value = 44400
$ws.Cells.Item(174, 9).Value = "Primera"
$ws.Cells.Item(174, 10).Value = 300
$ws.Cells.Item(174, 11).Value = 600
$ws.Cells.Item(174, 12).Value = 650
$ws.Cells.Item(174, 13).Value = 625
$ws.Cells.Item(174, 15).Value = "Provincia de Diguillín"
$ws.Cells.Item(174, 16).Value = 625

# Row 175
$ws.Cells.Item(175, 4).Value = 44181
$ws.Cells.Item(175, 9).Value = "Primera"
$ws.Cells.Item(175, 10).Value = 240
$ws.Cells.Item(175, 11).Value = 700
$ws.Cells.Item(175, 12).Value = 750
$ws.Cells.Item(175, 13).Value = 725
$ws.Cells.Item(175, 15).Value = "Provincia de Diguillín"
$ws.Cells.Item(175, 16).Value = 725

# Row 176
$ws.Cells.Item(176, 4).Value = 44494
$ws.Cells.Item(176, 9).Value = "Primera"
$ws.Cells.Item(176, 10).Value = 200
$ws.Cells.Item(176, 11).Value = 750
$ws.Cells.Item(176, 12).Value = 800
$ws.Cells.Item(176, 13).Value = 775
$ws.Cells.Item(176, 15).Value = "Región Metropolitana"
$ws.Cells.Item(176, 16).Value = 775

# Row 177
$ws.Cells.Item(177, 4).Value = 44342
$ws.Cells.Item(177, 9).Value = "Primera"
$ws.Cells.Item(177, 10).Value = 300
$ws.Cells.Item(177, 11).Value = 700
$ws.Cells.Item(177, 12).Value = 750
$ws.Cells.Item(177, 13).Value = 725
$ws.Cells.Item(177, 15).Value = "Región del Maule"
$ws.Cells.Item(177, 16).Value = 725

# Row 178
$ws.Cells.Item(178, 4).Value = 44328
$ws.Cells.Item(178, 9).Value = "Primera"
$ws.Cells.Item(178, 10).Value = 300
$ws.Cells.Item(178, 11).Value = 700
$ws.Cells.Item(178, 12).Value = 750
$ws.Cells.Item(178, 13).Value = 725
$ws.Cells.Item(178, 15).Value = "Región del Maule"
$ws.Cells.Item(178, 16).Value = 725

# Row 179
$ws.Cells.Item(179, 4).Value = 44301
$ws.Cells.Item(179, 9).Value = "Primera"
$ws.Cells.Item(179, 10).Value = 300
$ws.Cells.Item(179, 11).Value = 700
$ws.Cells.Item(179, 12).Value = 750
$ws.Cells.Item(179, 13).Value = 725
$ws.Cells.Item(179, 15).Value = "Provincia de Diguillín"
$ws.Cells.Item(179, 16).Value = 725

# Row 180
$ws.Cells.Item(180, 4).Value = 44301
$ws.Cells.Item(180, 9).Value = "Segunda"
$ws.Cells.Item(180, 10).Value = 80
$ws.Cells.Item(180, 11).Value = 600
$ws.Cells.Item(180, 12).Value = 600
$ws.Cells.Item(180, 13).Value = 600
$ws.Cells.Item(180, 15).Value = "Provincia de Diguillín"
$ws.Cells.Item(180, 16).Value = 600

# Row 181
$ws.Cells.Item(181, 4).Value = 44330
$ws.Cells.Item(181, 9).Value = "Primera"
$ws.Cells.Item(181, 10).Value = 300
$ws.Cells.Item(181, 11).Value = 700
$ws.Cells.Item(181, 12).Value = 750
$ws.Cells.Item(181, 13).Value = 725
$ws.Cells.Item(181, 15).Value = "Región del Maule"
$ws.Cells.Item(181, 16).Value = 725

# Row 182
$ws.Cells.Item(182, 4).Value = 44509
$ws.Cells.Item(182, 9).Value = "Primera"
$ws.Cells.Item(182, 10).Value = 300
$ws.Cells.Item(182, 11).Value = 650
$ws.Cells.Item(182, 12).Value = 700
$ws.Cells.Item(182, 13).Value = 675
$ws.Cells.Item(182, 15).Value = "Región del Maule"
$ws.Cells.Item(182, 16).Value = 675

# Row 183
$ws.Cells.Item(183, 4).Value = 44421
$ws.Cells.Item(183, 9).Value = "Primera"
$ws.Cells.Item(183, 10).Value = 300
$ws.Cells.Item(183, 11).Value = 650
$ws.Cells.Item(183, 12).Value = 700
$ws.Cells.Item(183, 13).Value = 675
$ws.Cells.Item(183, 15).Value = "Provincia de Diguillín"
$ws.Cells.Item(183, 16).Value = 675

# Row 184
$ws.Cells.Item(184, 4).Value = 44244
$ws.Cells.Item(184, 9).Value = "Primera"
$ws.Cells.Item(184, 10).Value = 300
$ws.Cells.Item(184, 11).Value = 700
$ws.Cells.Item(184, 12).Value = 750
$ws.Cells.Item(184, 13).Value = 725
$ws.Cells.Item(184, 15).Value = "Región del Maule"
$ws.Cells.Item(184, 16).Value = 725

# Row 185
$ws.Cells.Item(185, 4).Value = 44307
$ws.Cells.Item(185, 9).Value = "Primera"
$ws.Cells.Item(185, 10).Value = 300
$ws.Cells.Item(185, 11).Value = 600
$ws.Cells.Item(185, 12).Value = 650
$ws.Cells.Item(185, 13).Value = 625
$ws.Cells.Item(185, 15).Value = "Región del Maule"
$ws.Cells.Item(185, 16).Value = 625

# Row 186
$ws.Cells.Item(186, 4).Value = 44433
$ws.Cells.Item(186, 9).Value = "Primera"
$ws.Cells.Item(186, 10).Value = 300
$ws.Cells.Item(186, 11).Value = 700
$ws.Cells.Item(186, 12).Value = 750
$ws.Cells.Item(186, 13).Value = 725
$ws.Cells.Item(186, 15).Value = "Provincia de Diguillín"
$ws.Cells.Item(186, 16).Value = 725

# Row 187
$ws.Cells.Item(187, 4).Value = 44302
$ws.Cells.Item(187, 9).Value = "Primera"
$ws.Cells.Item(187, 10).Value = 300
$ws.Cells.Item(187, 11).Value = 600
$ws.Cells.Item(187, 12).Value = 650
$ws.Cells.Item(187, 13).Value = 625
$ws.Cells.Item(187, 15).Value = "Región del Maule"
$ws.Cells.Item(187, 16).Value = 625

# Row 188
$ws.Cells.Item(188, 4).Value = 44179
$ws.Cells.Item(188, 9).Value = "Primera"
$ws.Cells.Item(188, 10).Value = 160
$ws.Cells.Item(188, 11).Value = 700
$ws.Cells.Item(188, 12).Value = 750
$ws.Cells.Item(188, 13).Value = 725
$ws.Cells.Item(188, 15).Value = "Provincia de Diguillín"
$ws.Cells.Item(188, 16).Value = 725

# Row 189
$ws.Cells.Item(189, 4).Value = 44179
$ws.Cells.Item(189, 9).Value = "Segunda"
$ws.Cells.Item(189, 10).Value = 120
$ws.Cells.Item(189, 11).Value = 600
$ws.Cells.Item(189, 12).Value = 600
$ws.Cells.Item(189, 13).Value = 600
$ws.Cells.Item(189, 15).Value = "Provincia de Diguillín"
$ws.Cells.Item(189, 16).Value = 600

# New row 190 (sheet now spans A1:R190)
$ws.Cells.Item(190, 1).Value = 7
$ws.Cells.Item(190, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(190, 3).Value = "Ñuble"
$ws.Cells.Item(190, 4).Value = 44491
$ws.Cells.Item(190, 5).Value = 16
$ws.Cells.Item(190, 6).Value = 100112023
$ws.Cells.Item(190, 7).Value = "Brócoli"
$ws.Cells.Item(190, 8).Value = "Sin especificar"
$ws.Cells.Item(190, 9).Value = "Primera"
$ws.Cells.Item(190, 10).Value = 240
$ws.Cells.Item(190, 11).Value = 750
$ws.Cells.Item(190, 12).Value = 800
$ws.Cells.Item(190, 13).Value = 775
$ws.Cells.Item(190, 14).Value = "$/unidad"
$ws.Cells.Item(190, 15).Value = "Región Metropolitana"
$ws.Cells.Item(190, 16).Value = 775
$ws.Cells.Item(190, 17).Value = 1
$ws.Cells.Item(190, 18).Value = "Hortaliza"
$ws.Cells.Item(190, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
